$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-169 (11 -> 21)
$ws.Range("A2:A169").Value = 21

# Rows 170-337 (22 -> 20.25)
$ws.Range("A170:A337").Value = 20.25

# Rows 338-505 (19 -> 10.75)
$ws.Range("A338:A505").Value = 10.75
